# close #96 people met system
# Add a new "met" (encounter) row to the SceneQuest table on Sheet1 and
# keep the table/selection/dimension metadata in sync, mirroring what a
# user would do interactively in Excel after typing a new row into the
# table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New row of quest data (row 28): a "chance encounter" (偶遇) quest that
# uses the shared "met" identifier for its Ename/Figue/Script columns.
$ws.Range("A28").Value = 42010018
$ws.Range("B28").Value = "偶遇"
$ws.Range("C28").Value = 0
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = "met"
$ws.Range("F28").Value = "met"
$ws.Range("G28").Value = "met"
$ws.Range("I28").Value = 1
$ws.Range("J28").Value = "oneline"
$ws.Range("O28").Value = 100
$ws.Range("R28").Value = 1
$ws.Range("U28").Value = 200
$ws.Range("V28").Value = 200

# Grow the worksheet table ("表3") so the new row is included.
$tbl = $ws.ListObjects.Item(1)
[void]$tbl.Resize($ws.Range("A3:Z28"))

# Match the selection left behind after entering the new row.
[void]$ws.Range("F27").Select()
